# Update for insert release-notes.md f80ed2bb9e1dd81abc71d13817b8a44a756cee80
#
# Metadata sheet: bump version/status/date, add real contact display.
$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "0.4.0-snapshot-1"               # Version
$meta.Range("B6").Value  = "draft"                            # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"        # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"     # Contact

# Elements sheet: a new "business mapping" column was inserted before the
# existing RIM Mapping column upstream, which - in this flattened two-column
# view - shows up as the AK/AL mapping columns swapping places (header text,
# column width and the per-row mapping values all move together).
$elements = $wb.Worksheets.Item("Elements")

$lastRow = 6
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)   # column AK
    $alCell = $elements.Cells.Item($r, 38)   # column AL

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Column widths follow the content: AK becomes the wide "Spécification
# métier" column, AL becomes the narrower "RIM Mapping" column.
$elements.Columns.Item(37).ColumnWidth = 70.15
$elements.Columns.Item(38).ColumnWidth = 24.15
